$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM-derived values (re-run of scripts with new TPM data).
# Each row corresponds to a Sending/Target cluster combination for the
# Lama2 (ligand) - Rpsa (receptor) pair. Columns follow the header row:
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I/J Ligand derived specificity (avg/total),
# M Receptor average expression value, N Receptor total expression value,
# O/P Receptor derived specificity (avg/total),
# Q/R Edge average/total expression weight,
# S/T Edge average/total expression derived specificity.

$rows = @(
    @{ Row = 2;  E = 3; F = 1; G = 1.069299666666667;  H = 3.207899;  I = 0.003616700200628781; J = 0.003616700200628781; M = 13.929953;        N = 41.789859;  O = 0.09674275490334808; P = 0.09674275490334808; Q = 14.89529409958233;  R = 134.057646896241;  S = 0.00034988954106832;  T = 0.00034988954106832 }
    @{ Row = 3;  E = 3; F = 1; G = 1.069299666666667;  H = 3.207899;  I = 0.003616700200628781; J = 0.003616700200628781; M = 81.07766966666667;                         O = 0.5630799418129374; P = 0.5630799418129373; Q = 86.69632514867678;  R = 780.266926338091;  S = 0.002036491338524893; T = 0.002036491338524892 }
    @{ Row = 4;  E = 3; F = 1; G = 1.069299666666667;  H = 3.207899;  I = 0.003616700200628781; J = 0.003616700200628781; M = 48.98200233333333;  N = 146.946007; O = 0.3401773032837146;  P = 0.3401773032837146;  Q = 52.37643876769923;  R = 471.3879489092931; S = 0.001230319321035568;  T = 0.001230319321035568 }
    @{ Row = 5;                                          I = 0.8238194745364892;  J = 0.8238194745364891; M = 13.929953;        N = 41.789859;  O = 0.09674275490334808; P = 0.09674275490334808; Q = 3392.880990260406;  R = 30535.92891234365; S = 0.07969856550968858;  T = 0.07969856550968857 }
    @{ Row = 6;                                          I = 0.8238194745364892;  J = 0.8238194745364891; M = 81.07766966666667;                         O = 0.5630799418129374; P = 0.5630799418129373; Q = 19747.86879371711;                          S = 0.463876221786371;    T = 0.4638762217863708 }
    @{ Row = 7;                                          I = 0.8238194745364892;  J = 0.8238194745364891; M = 48.98200233333333;  N = 146.946007; O = 0.3401773032837146;  P = 0.3401773032837146;  Q = 11930.41387732303;  R = 107373.7248959072; S = 0.2802446872404297;   T = 0.2802446872404297 }
    @{ Row = 8;  G = 51.01955666666666; H = 153.05867;  I = 0.1725638252628821; J = 0.1725638252628821;                             M = 13.929953;        N = 41.789859;  O = 0.09674275490334808; P = 0.09674275490334808; Q = 710.7000264475032;  R = 6396.300238027529; S = 0.01669429985259119;  T = 0.01669429985259119 }
    @{ Row = 9;  G = 51.01955666666666; H = 153.05867;  I = 0.1725638252628821; J = 0.1725638252628821;                             M = 81.07766966666667;                         O = 0.5630799418129374; P = 0.5630799418129373; Q = 4136.54676195978;   R = 37228.92085763802; S = 0.09716722868804153;  T = 0.09716722868804152 }
    @{ Row = 10; G = 51.01955666666666; H = 153.05867;  I = 0.1725638252628821; J = 0.1725638252628821;                             M = 48.98200233333333;  N = 146.946007; O = 0.3401773032837146;  P = 0.3401773032837146;  Q = 2499.040043692299;  R = 22491.36039323069; S = 0.05870229672224938;  T = 0.05870229672224938 }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row
    foreach ($col in @('E','F','G','H','I','J','M','N','O','P','Q','R','S','T')) {
        if ($rowData.ContainsKey($col)) {
            $ws.Range("$col$r").Value = $rowData[$col]
        }
    }
}
